# Update countries & provincias Spain
# - Reorder two pairs of country rows (Suiza/Portugal and Mali/Albania/Guinea Ecuatorial)
#   by swapping their displayed country names (row data follows the name).
# - Refresh the "Datos actualizados" timestamp.
# - Update several numeric stats cells with newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 14:35"

# --- Swap Suiza (row 30) / Portugal (row 31) so Portugal now appears first ---
$ws.Range("A30").Value = "Portugal"
$ws.Range("A31").Value = "Suiza"

# --- Rotate Mali (row 110) / Albania (row 111) / Guinea Ecuatorial (row 112)
#     so Guinea Ecuatorial now appears first, then Mali, then Albania ---
$ws.Range("A110").Value = "Guinea Ecuatorial"
$ws.Range("A111").Value = "Mali"
$ws.Range("A112").Value = "Albania"

# --- Numeric updates ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1686791
$ws.Range("C4").Value = 355
$ws.Range("E4").Value = 1135735
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 99311

# Row 24: Paises Bajos
$ws.Range("B24").Value = 45445
$ws.Range("C24").Value = 209
$ws.Range("G24").Value = 8
$ws.Range("H24").Value = 5830

# Row 30: now Portugal
$ws.Range("B30").Value = 30788
$ws.Range("C30").Value = 165
$ws.Range("D30").Value = 17822
$ws.Range("E30").Value = 11636
$ws.Range("G30").Value = 14
$ws.Range("H30").Value = 1330

# Row 31: now Suiza
$ws.Range("B31").Value = 30746
$ws.Range("C31").Value = 10
$ws.Range("D31").Value = 28100
$ws.Range("E31").Value = 739
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 1907

# Row 55: Kazajistan
$ws.Range("D55").Value = 4506
$ws.Range("E55").Value = 3990

# Row 70: Azerbaiyan
$ws.Range("B70").Value = 4271
$ws.Range("C70").Value = 149
$ws.Range("D70").Value = 2741
$ws.Range("E70").Value = 1479
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 51

# Row 86: Croacia
$ws.Range("D86").Value = 2035
$ws.Range("E86").Value = 109
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 100

# Row 110: now Guinea Ecuatorial
$ws.Range("B110").Value = 1043
$ws.Range("C110").Value = 83
$ws.Range("D110").Value = 165
$ws.Range("E110").Value = 866
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 12

# Row 111: now Mali
$ws.Range("B111").Value = 1030
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 597
$ws.Range("E111").Value = 368
$ws.Range("H111").Value = 65

# Row 112: now Albania
$ws.Range("B112").Value = 1004
$ws.Range("C112").Value = 6
$ws.Range("D112").Value = 795
$ws.Range("E112").Value = 177
$ws.Range("H112").Value = 32

# Row 128: San Marino
$ws.Range("B128").Value = 666
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 270
$ws.Range("E128").Value = 354

# Row 158: Benin
$ws.Range("D158").Value = 84
$ws.Range("E158").Value = 104
